$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1720
$ws1.Range("F6").Value = 623
$ws1.Range("F7").Value = 1127
$ws1.Range("F8").Value = 1563
$ws1.Range("F10").Value = 12
$ws1.Range("F11").Value = 1481
$ws1.Range("F12").Value = 3111
$ws1.Range("F13").Value = 657
$ws1.Range("F14").Value = 1797
$ws1.Range("F15").Value = 1810
$ws1.Range("F16").Value = 863
$ws1.Range("F17").Value = 283
$ws1.Range("F18").Value = 3
$ws1.Range("F19").Value = 1489
$ws1.Range("F20").Value = 292
$ws1.Range("F22").Value = 18
$ws1.Range("F23").Value = 1245
$ws1.Range("F24").Value = 416
$ws1.Range("F25").Value = 472
$ws1.Range("F26").Value = 146
$ws1.Range("F27").Value = 4866
$ws1.Range("F28").Value = 2651
$ws1.Range("F29").Value = 760
$ws1.Range("F31").Value = 1672
$ws1.Range("F32").Value = 77
$ws1.Range("F33").Value = 181

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 32
$ws2.Range("F4").Value = 85
$ws2.Range("G4").Value = "不可售"
$ws2.Range("F5").Value = 25

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 44

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 44
$ws4.Range("F5").Value = 32
$ws4.Range("F7").Value = 85
$ws4.Range("G7").Value = "不可售"
$ws4.Range("F8").Value = 25
$ws4.Range("F9").Value = 1720
$ws4.Range("F11").Value = 623
$ws4.Range("F12").Value = 1127
$ws4.Range("F13").Value = 1563
$ws4.Range("F17").Value = 12
$ws4.Range("F18").Value = 1481
$ws4.Range("F19").Value = 3111
$ws4.Range("F20").Value = 657
$ws4.Range("F21").Value = 1797
$ws4.Range("F22").Value = 1810
$ws4.Range("F23").Value = 863
$ws4.Range("F24").Value = 283
$ws4.Range("F25").Value = 3
$ws4.Range("F26").Value = 1489
$ws4.Range("F27").Value = 292
$ws4.Range("F30").Value = 18
$ws4.Range("F32").Value = 1245
$ws4.Range("F33").Value = 416
$ws4.Range("F34").Value = 472
$ws4.Range("F35").Value = 146
$ws4.Range("F36").Value = 4866
$ws4.Range("F37").Value = 2655
$ws4.Range("F38").Value = 760
$ws4.Range("F40").Value = 1672
$ws4.Range("F43").Value = 77
$ws4.Range("F44").Value = 181
